{"js": "// The substabelecimento.docx template used the generic merge-field\n// names `cidade` and `data` for the place/date of signature line.\n// They are renamed to the more specific `cidade_assinatura` and\n// `data_assinatura` (see commit message: \"Vari\u00e1veis data de\n// assinatura e cidade de assinatura corrigidas no docx\").\n//\n// Before: {{ cidade }}, {{ data }}.\n// After : {{ cidade_assinatura }}, {{ data_assinatura }}.\nconst body = context.document.body;\n\nconst target = \"{{ cidade }}, {{ data }}.\";\nconst replacement = \"{{ cidade_assinatura }}, {{ data_assinatura }}.\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole phrase in one shot so \"cidade\" / \"data\" text\n  // elsewhere in the document (there is none, but just in case) is\n  // left untouched.\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n} else {\n  // Fallback: do the two renames individually, scoped to the exact\n  // merge-field tokens \"{{ cidade }}\" and \"{{ data }}\" so no other\n  // occurrence of those words is affected.\n  const cidadeResults = body.search(\"{{ cidade }}\", { matchCase: true });\n  cidadeResults.load(\"text\");\n  await context.sync();\n  for (const r of cidadeResults.items) {\n    r.insertText(\"{{ cidade_assinatura }}\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n\n  const dataResults = body.search(\"{{ data }}\", { matchCase: true });\n  dataResults.load(\"text\");\n  await context.sync();\n  for (const r of dataResults.items) {\n    r.insertText(\"{{ data_assinatura }}\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The substabelecimento.docx template used the generic merge-field\n# names `cidade` and `data` for the place/date of signature line.\n# They are renamed to the more specific `cidade_assinatura` and\n# `data_assinatura` (see commit message: \"Vari\u00e1veis data de\n# assinatura e cidade de assinatura corrigidas no docx\").\n#\n# Before: {{ cidade }}, {{ data }}.\n# After : {{ cidade_assinatura }}, {{ data_assinatura }}.\n$d = $word.ActiveDocument\n\n$target = \"{{ cidade }}, {{ data }}.\"\n$replacement = \"{{ cidade_assinatura }}, {{ data_assinatura }}.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$found = $find.Execute(\n    $target,       # FindText\n    $false,        # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap (wdFindContinue)\n    $false,        # Format\n    $replacement,  # ReplaceWith\n    2              # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    # Fallback: rename the two merge fields individually, scoped to\n    # the exact tokens \"{{ cidade }}\" / \"{{ data }}\" so no other\n    # occurrence of those words elsewhere in the document is touched.\n    $findCidade = $d.Content.Find\n    $findCidade.ClearFormatting()\n    $findCidade.Replacement.ClearFormatting()\n    $findCidade.Execute(\n        \"{{ cidade }}\", $false, $false, $false, $false, $false,\n        $true, 1, $false, \"{{ cidade_assinatura }}\", 2\n    )\n\n    $findData = $d.Content.Find\n    $findData.ClearFormatting()\n    $findData.Replacement.ClearFormatting()\n    $findData.Execute(\n        \"{{ data }}\", $false, $false, $false, $false, $false,\n        $true, 1, $false, \"{{ data_assinatura }}\", 2\n    )\n}\n"}
